# Conserto do erro com o rotulo da coluna 2050 nas tabelas
# e retirada das linhas com total das tabelas.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Fix the mis-typed "2050" (or "2041-2050") header label in column E.
#    The cell currently holds a stray numeric value
#    (673.701694546605) instead of the intended text label.
#    We write the label as a formula first, then paste its computed
#    result back as a value so the cell ends up holding literal text
#    (not a number, and not a live formula) while keeping its existing
#    style untouched.
# ---------------------------------------------------------------------
$sheetLabels = @{
    1 = "2050"
    2 = "2050"
    3 = "2050"
    4 = "2041-2050"
    5 = "2050"
}

foreach ($idx in $sheetLabels.Keys) {
    $ws = $wb.Worksheets.Item($idx)
    $cell = $ws.Cells.Item(1, 5)
    $cell.Formula = "=""" + $sheetLabels[$idx] + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

# ---------------------------------------------------------------------
# 2) Remove the "Total" row from each table.
#    Sheets 1-4 have it at row 13, sheet 6 at row 4.
#    Sheet 5 never had a Total row.
# ---------------------------------------------------------------------
foreach ($idx in 1, 2, 3, 4) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Rows.Item(13).Delete()
}

$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows.Item(4).Delete()
